$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 19:26 which no longer exist in the updated data range
$ws.Rows("19:26").Delete()

# New header for column F, matching style (bold, border, centered) of existing headers
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "weights"

# Updated data values for A2:F18
$data = @(
    @(0.292, 0.8284900000000001, 0.292, 0.7905203322205661, 0.7905203322205661, 1),
    @(0.314, 0.8284900000000001, 0.02200000000000002, 0.7905203322205661, 0.7905203322205661, 0.07006369426751599),
    @(0.341, 0.8284900000000001, 0.02700000000000002, 0.7905203322205661, 0.7905203322205659, 0.0791788856304986),
    @(0.367, 0.8284900000000001, 0.02599999999999997, 0.7905203322205661, 0.7905203322205659, 0.07084468664850127),
    @(0.398, 0.8284900000000001, 0.03100000000000003, 0.7905203322205661, 0.7905203322205659, 0.07788944723618096),
    @(0.443, 0.8284900000000001, 0.04499999999999998, 0.7905203322205661, 0.7905203322205659, 0.1015801354401805),
    @(0.489, 0.81699, 0.04599999999999999, 0.779547376819129, 0.6738730454531172, 0.09406952965235171),
    @(0.542, 0.81699, 0.05300000000000005, 0.779547376819129, 0.6738730454531172, 0.09778597785977867),
    @(0.587, 0.80548, 0.04499999999999993, 0.768564879717343, 0.6738730454531172, 0.07666098807495729),
    @(0.64, 0.7939700000000001, 0.05300000000000005, 0.757582382615557, 0.6738730454531172, 0.08281250000000007),
    @(0.708, 0.78247, 0.06799999999999995, 0.7466094272141199, 0.6738730454531172, 0.09604519774011293),
    @(0.765, 0.7709600000000001, 0.05700000000000005, 0.735626930112334, 0.6738730454531172, 0.0745098039215687),
    @(0.8179999999999999, 0.75945, 0.05299999999999994, 0.7246444330105478, 0.6738730454531172, 0.06479217603911973),
    @(0.867, 0.74795, 0.04900000000000004, 0.713671477609111, 0.6738730454531172, 0.05651672433679359),
    @(0.924, 0.7364400000000001, 0.05700000000000005, 0.702688980507325, 0.6738730454531172, 0.06168831168831174),
    @(0.962, 0.72493, 0.03799999999999992, 0.6917064834055389, 0.6738730454531172, 0.03950103950103942),
    @(0.991, 0.715, 0.02900000000000003, 0.6822315749589067, 0.6738730454531172, 0.029263370332997)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $row++
}
